$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.787.87'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '2.080.53'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Formula = '''233.17'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Formula = '''0.626'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').Formula = '''58.65'
$ws.Range('E7').Value = '  +2.42%  '
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').Formula = '''0.0789'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('E11').Value = '  +4.20%  '
$ws.Range('D12').Value = '2.387.50'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').Formula = '''14.86'
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').Formula = '''21.26'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Formula = '''0.779'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('D17').Value = '2.065.03'
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').Value = '37.765.53'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').Formula = '''71.69'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Formula = '''6.12'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('D21').Value = '0.0₃0846'
$ws.Range('E21').Value = '  +2.61%  '
$ws.Range('D22').Formula = '''228.69'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Formula = '''2.41'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Formula = '''9.59'
$ws.Range('E26').Value = '  +7.30%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Formula = '''171.26'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').Formula = '''0.139'
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('D29').Formula = '''1.41'
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('D30').Formula = '''19.51'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').Formula = '''0.122'
$ws.Range('E31').Value = '  +2.30%  '
$ws.Range('E32').Value = '  +2.77%  '
$ws.Range('D33').Formula = '''0.0634'
$ws.Range('E33').Value = '  +1.14%  '
$ws.Range('E34').Value = '  +1.35%  '
$ws.Range('D35').Formula = '''2.47'
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('D36').Formula = '''3.44'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').Formula = '''1.83'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').Formula = '''5.51'
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Formula = '''0.0235'
$ws.Range('E40').Value = '  +9.73%  '
$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D41').Formula = '''0.0975'
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Formula = '''99.62'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').Formula = '''2.94'
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Formula = '''17.21'
$ws.Range('E44').Value = '  +9.62%  '
$ws.Range('D45').Value = '1.452.35'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Formula = '''1.07'
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').Formula = '''4.14'
$ws.Range('E48').Value = '  -3.60%  '
$ws.Range('D49').Formula = '''7.39'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').Formula = '''3.00'
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('D51').Value = '2.273.08'
$ws.Range('E51').Value = '  +0.18%  '
